$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 72 data (mirrors row 71's layout/style)
$ws.Range("A71:F71").Copy()
$ws.Range("A72:F72").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Cells.Item(72, 1).Value = 43791
$ws.Cells.Item(72, 2).Value = 2189.3936572804901
$ws.Cells.Item(72, 3).Value = 2207.0300000000002
$ws.Range("D72").FormulaR1C1 = "=100*(RC[-2]-RC[-1])/RC[-1]"
$ws.Cells.Item(72, 5).Value = 169
$ws.Cells.Item(72, 6).Value = "Crm opened 11/19/2019"

# Update selection
$ws.Range("E73").Select()
